# "Generate Report for Handback" - refresh the localization-status report
# after a handback has completed for 2c29110a-a729-4e33-bcb9-89a01b98271b:
#   - status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#   - the "Latest Handback DateTime" columns get the new handback timestamps
#   - the stale "handback file is not latest" error note is cleared
#   - column P (Error Detail) narrows back down now that it holds no long text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is 2c29110a-a729-4e33-bcb9-89a01b98271b.md ---
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is 2c29110a-a729-4e33-bcb9-89a01b98271b.md ---
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-23 14:58:49"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 13

# --- de-de sheet: row 3 is 2c29110a-a729-4e33-bcb9-89a01b98271b.md ---
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-23 14:58:57"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 13

Write-Output "Report regenerated for handback."
